$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 374, shifting existing rows 374-396 down to 375-397.
$ws.Rows.Item(374).Insert()

# Populate the newly inserted row 374 with the new weekly price observation.
$ws.Range("A374").Value = 9
$ws.Range("B374").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C374").Value = "Metropolitana"
$ws.Range("D374").Value = 44746
$ws.Range("E374").Value = 13
$ws.Range("F374").Value = 100112032
$ws.Range("G374").Value = "Zapallo italiano"
$ws.Range("H374").Value = "Sin especificar"
$ws.Range("I374").Value = "Primera"
$ws.Range("J374").Value = 97
$ws.Range("K374").Value = 11000
$ws.Range("L374").Value = 13000
$ws.Range("M374").Value = 11990
$ws.Range("N374").Value = "$/caja 50 unidades"
$ws.Range("O374").Value = "Región de Arica y Parinacota"
$ws.Range("P374").Value = 240
$ws.Range("Q374").Value = 50
$ws.Range("R374").Value = "Hortaliza"
